# Apply Resident Services requirements edits to the 'Details' sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")

$ws.Range("O2").Value2 = 'Research info '
$ws.Range("N3").Value2 = '1. what is the use case of locking specific Biometric auths and combinations?' + [char]10 + '2. What is the data taken as input from the Individual?' + [char]10 + '3. Is there a mechanism to lock OTP Authentication?' + [char]10 + '4. need more clarity on a2' + [char]10 + '5. Can Lock/Unlock will be perfomed only by OTP authentication of Mobile number or can it be done by email as well?' + [char]10 + '6. Need to check with IDA on the process of authentication done, so that the gaps could be covered' + [char]10 + '7. Can both UIN and VID be locked at the same time?' + [char]10 + 'Is this feature exclusive to Resident services?' + [char]10 + 'what is Demographic auth?' + [char]10 + '*Notification service is required here' + [char]10 + ''
$ws.Range("O3").Value2 = 'Resident can lock his UID number Via Resident Portal, by doing this Resident cannot perform any sort of Authentication by using UID, UID Token & ANCS Token for Biometric, Demographic & OTP Based authentication. Once UID is locked resident can authenticate using 16 digit VID Number for all forms of authentication (Demo, Bio & OTP).' + [char]10 + ''
$ws.Range("N4").Value2 = '1. What is the use case of reprinting? Will there be a cost associated to it? If not can it be abused by the individual?will there be a limit on number times an individual can access it? If cost is associated, will there be a check performed for the payment?' + [char]10 + '2. Why RID is accepted as an input parameter? What is the use case.' + [char]10 + '3. Aadhar provides option to enter non registered mobile number on which the OTP can be received.' + [char]10 + 'Is this feature part of Reg Proc(though not Re-print but Print of UIN)?' + [char]10 + '* Notification service is required'
$ws.Range("O4").Value2 = 'Aadhar asks for Payment to perform Re-print' + [char]10 + '"Order Aadhaar Reprint" is a new service launched by UIDAI w.e.f. 01-12-2018 on Pilot basis which facilitates the residents of India to get their Aadhaar letter reprint by paying nominal charges, in case, Aadhaar letter of resident has been lost, misplaced or if they want a new copy. Residents who do not have registered mobile number can also "Order Aadhaar Reprint" using Non-Registered / Alternate Mobile Number.'
$ws.Range("O7").Value2 = 'In addition to English you can update/do correction in your address in any of the following languages:' + [char]10 + 'Assamese, Bengali, English, Gujarati, Hindi, Kannada, Malayalam, Marathi, Odia, Punjabi, Tamil, Telugu and Urdu.' + [char]10 + 'Q:I do not have any document proof of my address. Can I still update my address in my Aadhaar? Address Verifier concept' + [char]10 + 'Q:Does submission of request guarantee Updation of information?' + [char]10 + 'Submission of information for update does not guarantee update of Aadhaar data. The information submitted is subject to verification and validation. Furnishing of incorrect information/suppression of information would lead to rejection of application.'
$ws.Range("N8").Value2 = '1. why is RID an input here and not UIN?' + [char]10 + '2. What are the Statuses required?'
$ws.Range("N8").WrapText = $true
$ws.Range("N9").Value2 = '1.what is security code?' + [char]10 + '2. What is Statuses required?'
$ws.Range("N9").WrapText = $true
$ws.Range("N10").Value2 = 'e-UIN should be password protected (Postal Code/combination of Name-DoB - TBD)?' + [char]10 + 'Should the RID also be considered for downloading e-UIN'
$ws.Range("O10").Value2 = 'Resident can download e-Aadhaar using 28 digit enrolment no. along with Full Name and Pin code. In this download process OTP is received on registered mobile no. Resident can also use TOTP to download e-Aadhaar instead of OTP. TOTP can generated using mAadhaar mobile Application.' + [char]10 + 'e-Aadhaar is a password protected electronic copy of Aadhaar, which is digitally signed by the competent Authority of UIDAI.'
$ws.Range("N11").Value2 = 'exact paramters needs to be known. For this the Authentication data stored by IDA needs to be known.' + [char]10 + 'Use case for an Individual to access the Auth history needs to be known, as it will help in understanding the Paramters to be displayed.'

# Restore focus/selection to the sheet and cell the author ended editing on
$ws.Activate()
$ws.Range("N4").Select()
